$wb = $excel.ActiveWorkbook

# "Test 1" sheet: fix shared string "Finland" -> "France" (cell A1)
$ws1 = $wb.Worksheets.Item("Test 1")
$ws1.Range("A1").Value = "France"

# Update the cached RPP values in column E, rows 18-22 on "Test 1".
# "Test 2" pulls these via formulas ('Test 1'!E18 etc.) and will recalc.
$ws1.Range("E18").Value = 1.7
$ws1.Range("E19").Value = 0.6
$ws1.Range("E20").Value = -0.1
$ws1.Range("E21").Value = 0.6
$ws1.Range("E22").Value = 1.2

# Move the active selection on "Test 1" from L19 to F22
$ws1.Range("F22").Select()
